# OW-535 - margin call generation fix
# Test fixture update: rename Portfolio ID value "p1" -> "p1a" on the
# OneBilateral.xlsx sample sheet (cell AP2, column "Portfolio ID"), and
# bring the saved view/selection state in line with where the author was
# working in the sheet when the fixture was re-saved (scrolled further
# right, with the active cell moved down to row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data fix -------------------------------------------------------
# Cell AP2 holds the "Portfolio ID" test value; update it from "p1" to "p1a".
$ws.Range("AP2").Value = "p1a"

# --- View / selection state -----------------------------------------------
# Move the active selection to AM14 (previously AM2) and scroll the window
# so column AI becomes the left-most visible column (previously AB).
$excel.ActiveWindow.ScrollColumn = 35
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AM14").Select() | Out-Null
